$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.352.12"
$ws.Range("E2").Value = "  +4.22%  "

$ws.Range("D3").Value = "3.640.06"
$ws.Range("E3").Value = "  +3.20%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "202.01"
$ws.Range("E5").Value = "  +11.60%  "

$ws.Range("D6").Value = "568.37"
$ws.Range("E6").Value = "  -1.75%  "

$ws.Range("D7").Value = "3.617.70"
$ws.Range("E7").Value = "  +2.94%  "

$ws.Range("D8").Value = "0.619"
$ws.Range("E8").Value = "  +2.42%  "

$ws.Range("E9").Value = "  -0.34%  "

$ws.Range("D10").Value = "0.679"
$ws.Range("E10").Value = "  +2.61%  "

$ws.Range("D11").Value = "58.25"
$ws.Range("E11").Value = "  +8.77%  "

$ws.Range("E12").Value = "  +7.94%  "

$ws.Range("D13").Value = "0.0000293"
$ws.Range("E13").Value = "  +16.72%  "

$ws.Range("D14").Value = "10.05"
$ws.Range("E14").Value = "  +3.55%  "

$ws.Range("D15").Value = "4.212.87"
$ws.Range("E15").Value = "  +2.98%  "

$ws.Range("D16").Value = "3.637.26"
$ws.Range("E16").Value = "  +3.02%  "

$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("D18").Value = "68.210.23"
$ws.Range("E18").Value = "  +4.17%  "

$ws.Range("D19").Value = "18.60"
$ws.Range("E19").Value = "  +2.57%  "

$ws.Range("D20").Value = "12.42"
$ws.Range("E20").Value = "  +3.05%  "

$ws.Range("E21").Value = "  +3.99%  "

$ws.Range("D22").Value = "401.92"
$ws.Range("E22").Value = "  +3.15%  "

$ws.Range("D23").Value = "13.12"
$ws.Range("E23").Value = "  +27.82%  "

$ws.Range("D24").Value = "4.21"
$ws.Range("E24").Value = "  -1.01%  "

$ws.Range("D25").Value = "85.68"
$ws.Range("E25").Value = "  +1.81%  "

$ws.Range("E26").Value = "  +3.93%  "

$ws.Range("E27").Value = "  +2.76%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "3.85"
$ws.Range("E28").Value = "  +9.43%  "

$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").Value = "6.12"
$ws.Range("E29").Value = "  +2.04%  "

$ws.Range("D30").Value = "8.12"
$ws.Range("E30").Value = "  +20.49%  "

$ws.Range("D31").Value = "9.16"
$ws.Range("E31").Value = "  +3.89%  "

$ws.Range("D32").Value = "31.87"
$ws.Range("E32").Value = "  +3.80%  "

$ws.Range("D33").Value = "689.90"
$ws.Range("E33").Value = "  +13.74%  "

$ws.Range("D34").Value = "12.24"
$ws.Range("E34").Value = "  +2.48%  "

$ws.Range("E35").Value = "  +4.92%  "

$ws.Range("D36").Value = "64.37"
$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("D37").Value = "42.84"
$ws.Range("E37").Value = "  +4.82%  "

$ws.Range("D38").Value = "0.427"
$ws.Range("E38").Value = "  +15.80%  "

$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").Value = "0.0₃0780"
$ws.Range("E40").Value = "  +5.26%  "

$ws.Range("E41").Value = "  +8.14%  "

$ws.Range("D42").Value = "3.260.34"
$ws.Range("E42").Value = "  +15.09%  "

$ws.Range("D43").Value = "3.14"
$ws.Range("E43").Value = "  +14.18%  "

$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  +17.73%  "

$ws.Range("D45").Value = "3.03"
$ws.Range("E45").Value = "  +37.36%  "

$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("D47").Value = "0.0419"
$ws.Range("E47").Value = "  +3.29%  "

$ws.Range("D48").Value = "2.75"
$ws.Range("E48").Value = "  +11.33%  "

$ws.Range("D49").Value = "8.91"
$ws.Range("E49").Value = "  +9.32%  "

$ws.Range("E50").Value = "  +2.03%  "

$ws.Range("D51").Value = "3.10"
$ws.Range("E51").Value = "  +6.30%  "
